# Apply updated crypto price/volume data scraped on Fri Aug 18 08:31:28 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Preserve the cell's existing style/format while forcing the assigned
    # string to be stored as literal text, even if it looks numeric
    # (e.g. "219.89"), matching the workbook's inline-string convention.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "26.598.74"
$ws.Range("E2").Value = "  -7.20%  "
$ws.Range("D3").Value = "1.694.83"
$ws.Range("E3").Value = "  -5.97%  "
$ws.Range("E4").Value = "  +0.29%  "
Set-TextValue $ws.Range("D5") "219.89"
$ws.Range("E5").Value = "  -5.24%  "
Set-TextValue $ws.Range("D6") "0.5094"
$ws.Range("E6").Value = "  -13.91%  "
$ws.Range("E7").Value = "  +0.21%  "
Set-TextValue $ws.Range("D8") "0.2652"
$ws.Range("E8").Value = "  -4.41%  "
Set-TextValue $ws.Range("D9") "22.16"
$ws.Range("E9").Value = "  -4.83%  "
Set-TextValue $ws.Range("D10") "0.06288"
$ws.Range("E10").Value = "  -7.68%  "
Set-TextValue $ws.Range("D11") "0.07376"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").Value = "1.694.25"
$ws.Range("E12").Value = "  -5.93%  "
Set-TextValue $ws.Range("D13") "4.510"
$ws.Range("E13").Value = "  -5.33%  "
Set-TextValue $ws.Range("D14") "0.5849"
$ws.Range("E14").Value = "  -6.03%  "
$ws.Range("D15").Value = "1.925.58"
$ws.Range("E15").Value = "  -5.91%  "
Set-TextValue $ws.Range("D16") "0.000008395"
$ws.Range("E16").Value = "  -8.60%  "
Set-TextValue $ws.Range("D17") "65.53"
$ws.Range("E17").Value = "  -13.36%  "
$ws.Range("D18").Value = "26.629.33"
$ws.Range("E18").Value = "  -7.04%  "
Set-TextValue $ws.Range("D19") "5.018"
$ws.Range("E19").Value = "  -8.08%  "
$ws.Range("E20").Value = "  +0.18%  "
$ws.Range("E21").Value = "  -4.31%  "
Set-TextValue $ws.Range("D22") "186.23"
$ws.Range("E22").Value = "  -11.54%  "
Set-TextValue $ws.Range("D23") "6.271"
$ws.Range("E23").Value = "  -8.18%  "
Set-TextValue $ws.Range("D24") "1.007"
$ws.Range("E24").Value = "  +0.27%  "
Set-TextValue $ws.Range("D25") "144.72"
$ws.Range("E25").Value = "  -5.96%  "
Set-TextValue $ws.Range("D26") "7.521"
$ws.Range("E26").Value = "  -4.16%  "
Set-TextValue $ws.Range("D27") "0.1155"
$ws.Range("E27").Value = "  -8.91%  "
Set-TextValue $ws.Range("D28") "15.64"
$ws.Range("E28").Value = "  -4.98%  "
Set-TextValue $ws.Range("D29") "1.337"
$ws.Range("E29").Value = "  -5.26%  "
Set-TextValue $ws.Range("D30") "0.05668"
Set-TextValue $ws.Range("D31") "1.344"
$ws.Range("E31").Value = "  -5.68%  "
Set-TextValue $ws.Range("D32") "3.520"
$ws.Range("E32").Value = "  -6.66%  "
Set-TextValue $ws.Range("D33") "3.495"
$ws.Range("E33").Value = "  -6.60%  "
Set-TextValue $ws.Range("D34") "1.651"
$ws.Range("E34").Value = "  -4.47%  "
Set-TextValue $ws.Range("D35") "1.022"
$ws.Range("E35").Value = "  -3.36%  "
Set-TextValue $ws.Range("D36") "0.6048"
$ws.Range("E36").Value = "  -5.78%  "
Set-TextValue $ws.Range("D37") "2.363"
$ws.Range("E37").Value = "  -5.51%  "
Set-TextValue $ws.Range("D38") "2.678"
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D39").Value = "1.100.63"
$ws.Range("E39").Value = "  -4.32%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D40") "0.01612"
$ws.Range("E40").Value = "  -4.62%  "
Set-TextValue $ws.Range("D41") "0.8657"
$ws.Range("E41").Value = "  -1.98%  "
Set-TextValue $ws.Range("D42") "5.860"
$ws.Range("E43").Value = "  -0.41%  "
Set-TextValue $ws.Range("D44") "99.48"
$ws.Range("E44").Value = "  -0.60%  "
$ws.Range("D45").Value = "1.854.00"
$ws.Range("E45").Value = "  -5.19%  "
Set-TextValue $ws.Range("D46") "0.00000000109"
$ws.Range("E46").Value = "  -3.91%  "
Set-TextValue $ws.Range("D47") "56.73"
$ws.Range("E47").Value = "  -6.12%  "
Set-TextValue $ws.Range("D48") "8.165"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("E49").Value = "  +0.51%  "
Set-TextValue $ws.Range("D50") "0.05247"
$ws.Range("E50").Value = "  -4.09%  "
Set-TextValue $ws.Range("D51") "0.4329"
$ws.Range("E51").Value = "  -3.28%  "
